$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.136
$ws.Range("E14").Value = 13.235
$ws.Range("E16").Value = 12.989
$ws.Range("E21").Value = 13.449
$ws.Range("E23").Value = 13.136
$ws.Range("E25").Value = 12.659
